# Scheduled runner refresh: update currentAveragePrice*/LevePrice*/LeveProfit*
# columns (H:N) for the affected leve rows on each crafting-job sheet, as
# pulled from the latest market-board snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1744.2
$ws.Range("I51").Value = 1423.3334
$ws.Range("J51").Value = 1881.7142
$ws.Range("K51").Value = 1423.3334
$ws.Range("L51").Value = 1881.7142
$ws.Range("M51").Value = -939.3334
$ws.Range("N51").Value = -2849.7142

$ws.Range("H112").Value = 1767.8334
$ws.Range("I112").Value = 2194.8
$ws.Range("J112").Value = 1462.8572
$ws.Range("K112").Value = 6584.400000000001
$ws.Range("L112").Value = 4388.571599999999
$ws.Range("M112").Value = -5476.400000000001
$ws.Range("N112").Value = -6604.571599999999

$ws.Range("H134").Value = 49450
$ws.Range("J134").Value = 49450
$ws.Range("L134").Value = 49450
$ws.Range("N134").Value = -59590

$ws.Range("H136").Value = 49500
$ws.Range("J136").Value = 49500
$ws.Range("L136").Value = 49500
$ws.Range("N136").Value = -59700

$ws.Range("H139").Value = 45714.75
$ws.Range("J139").Value = 45714.75
$ws.Range("L139").Value = 45714.75
$ws.Range("N139").Value = -55994.75

$ws.Range("H140").Value = 38000
$ws.Range("J140").Value = 38000
$ws.Range("L140").Value = 38000
$ws.Range("N140").Value = -48360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 750.25
$ws.Range("I4").Value = 900.5
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 900.5
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = -784.5
$ws.Range("N4").Value = -832

$ws.Range("H5").Value = 400
$ws.Range("I5").Value = 440
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 440
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -328
$ws.Range("N5").Value = -524

$ws.Range("H45").Value = 2007
$ws.Range("I45").Value = 2008.8572
$ws.Range("K45").Value = 2008.8572
$ws.Range("M45").Value = -1631.8572

$ws.Range("H63").Value = 2500
$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2500
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -1814

$ws.Range("H66").Value = 2500
$ws.Range("I66").Value = 2500
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 12500
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -9068

$ws.Range("H68").Value = 22000
$ws.Range("J68").Value = 22000
$ws.Range("L68").Value = 22000
$ws.Range("N68").Value = -23622

$ws.Range("H71").Value = 22000
$ws.Range("J71").Value = 22000
$ws.Range("L71").Value = 66000
$ws.Range("N71").Value = -74112

$ws.Range("H74").Value = 10110131
$ws.Range("I74").Value = 8182321.5
$ws.Range("J74").Value = 27781716
$ws.Range("K74").Value = 8182321.5
$ws.Range("L74").Value = 27781716
$ws.Range("M74").Value = -8181447.5
$ws.Range("N74").Value = -27783464

$ws.Range("H77").Value = 10110131
$ws.Range("I77").Value = 8182321.5
$ws.Range("J77").Value = 27781716
$ws.Range("K77").Value = 40911607.5
$ws.Range("L77").Value = 138908580
$ws.Range("M77").Value = -40907239.5
$ws.Range("N77").Value = -138917316

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 400
$ws.Range("I4").Value = 440
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 440
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -325
$ws.Range("N4").Value = -530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 349
$ws.Range("I22").Value = 352.2
$ws.Range("J22").Value = 345.8
$ws.Range("K22").Value = 352.2
$ws.Range("L22").Value = 345.8
$ws.Range("M22").Value = -2.199999999999989
$ws.Range("N22").Value = -1045.8

$ws.Range("H58").Value = 5628.143
$ws.Range("I58").Value = 2600
$ws.Range("J58").Value = 7310.4443
$ws.Range("K58").Value = 2600
$ws.Range("L58").Value = 7310.4443
$ws.Range("M58").Value = -2397
$ws.Range("N58").Value = -7716.4443

$ws.Range("H122").Value = 3031968.8
$ws.Range("I122").Value = 5129194.5
$ws.Range("J122").Value = 2642.7778
$ws.Range("K122").Value = 15387583.5
$ws.Range("L122").Value = 7928.3334
$ws.Range("M122").Value = -15385133.5
$ws.Range("N122").Value = -12828.3334

$ws.Range("H136").Value = 5628.143
$ws.Range("I136").Value = 2600
$ws.Range("J136").Value = 7310.4443
$ws.Range("K136").Value = 7800
$ws.Range("L136").Value = 21931.3329
$ws.Range("M136").Value = -5250
$ws.Range("N136").Value = -27031.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1752.9412
$ws.Range("I16").Value = 1400
$ws.Range("J16").Value = 1775
$ws.Range("K16").Value = 4200
$ws.Range("L16").Value = 5325
$ws.Range("M16").Value = -4027
$ws.Range("N16").Value = -5671

$ws.Range("H112").Value = 8101.125
$ws.Range("I112").Value = 3206.75
$ws.Range("J112").Value = 9080
$ws.Range("K112").Value = 9620.25
$ws.Range("L112").Value = 27240
$ws.Range("M112").Value = -8512.25
$ws.Range("N112").Value = -29456

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2196.2727
$ws.Range("I7").Value = 1551
$ws.Range("J7").Value = 2565
$ws.Range("K7").Value = 1551
$ws.Range("L7").Value = 2565
$ws.Range("M7").Value = -1439
$ws.Range("N7").Value = -2789

$ws.Range("H40").Value = 5062.2856
$ws.Range("I40").Value = 5113.0625
$ws.Range("J40").Value = 4899.8
$ws.Range("K40").Value = 5113.0625
$ws.Range("L40").Value = 4899.8
$ws.Range("M40").Value = -4977.0625
$ws.Range("N40").Value = -5171.8

$ws.Range("H122").Value = 3700.087
$ws.Range("I122").Value = 3056.0908
$ws.Range("J122").Value = 4290.4165
$ws.Range("K122").Value = 9168.2724
$ws.Range("L122").Value = 12871.2495
$ws.Range("M122").Value = -6718.2724
$ws.Range("N122").Value = -17771.2495

$ws.Range("H126").Value = 2196.2727
$ws.Range("I126").Value = 1551
$ws.Range("J126").Value = 2565
$ws.Range("K126").Value = 4653
$ws.Range("L126").Value = 7695
$ws.Range("M126").Value = -2183
$ws.Range("N126").Value = -12635

$ws.Range("H136").Value = 4314649
$ws.Range("I136").Value = 6255626
$ws.Range("J136").Value = 1367.2222
$ws.Range("K136").Value = 18766878
$ws.Range("L136").Value = 4101.6666
$ws.Range("M136").Value = -18764328
$ws.Range("N136").Value = -9201.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 505900
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 755850
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 755850
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -757098

$ws.Range("H65").Value = 505900
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 755850
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 3779250
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -3785490

$ws.Range("H122").Value = 111113630
$ws.Range("I122").Value = 142859660
$ws.Range("J122").Value = 2502.5
$ws.Range("K122").Value = 428578980
$ws.Range("L122").Value = 7507.5
$ws.Range("M122").Value = -428576530
$ws.Range("N122").Value = -12407.5
